$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new book data
$ws.Range("A2").Value = "Agile Methodology Master the art of Software Development:"
$ws.Range("B2").Value = "978-0-618-26030-0"
$ws.Range("C2").Value = "['Jason Roy', 'John Jacob Henry Rose']"
$ws.Range("D2").Value = "['The MIT Press']"

# Remove row 3 (the old second book entry)
$ws.Rows("3:3").Delete()
